# Insert a new weekly price record row above row 696 (Sandia / Femacal de La Calera),
# shifting all subsequent rows down by one (old row 696 becomes 697, ..., old 756 becomes 757).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 696 (pushes existing row 696 and below down one row)
$ws.Rows.Item(696).Insert()

# Populate the new row 696 with the new weekly record
$ws.Cells.Item(696, 1).Value  = 3
$ws.Cells.Item(696, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(696, 3).Value  = "Coquimbo"
$ws.Cells.Item(696, 4).Value  = 45265
$ws.Cells.Item(696, 5).Value  = 5
$ws.Cells.Item(696, 6).Value  = 100112028
$ws.Cells.Item(696, 7).Value  = "Sandia"
$ws.Cells.Item(696, 8).Value  = "Sin especificar"
$ws.Cells.Item(696, 9).Value  = "Primera"
$ws.Cells.Item(696, 10).Value = 160
$ws.Cells.Item(696, 11).Value = 800
$ws.Cells.Item(696, 12).Value = 800
$ws.Cells.Item(696, 13).Value = 800
$ws.Cells.Item(696, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(696, 15).Value = "Perú"
$ws.Cells.Item(696, 16).Value = 800
$ws.Cells.Item(696, 17).Value = 1
$ws.Cells.Item(696, 18).Value = "Hortaliza"

# Give the new date cell (column D) the same date/time number format used by the other rows
$ws.Cells.Item(696, 4).NumberFormat = $ws.Cells.Item(697, 4).NumberFormat
